$d = $word.ActiveDocument

# "Versi" + "on"  ->  single run "Version"
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Replacement.ClearFormatting()
$find1.Execute("Version", $true, $false, $false, $false, $false, $true, 1, $false, "Version", 2)

# " 2"  ->  " 1"   (stays inside its own run, doesn't cross the _GoBack bookmark)
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute("2", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2)

# Move the trailing "." from after the bookmark to before it, so the
# " 1" run becomes " 1." and the bookmark ends up at the end of the text.
$bm = $d.Bookmarks("_GoBack")
$bmRange = $bm.Range
$bmRange.InsertBefore(".")

# Remove the now-orphaned "." run that trails the (re-positioned) bookmark.
$bm2 = $d.Bookmarks("_GoBack")
$trailing = $bm2.Range
$trailing.MoveEnd(1, 1)
$trailing.Delete()
